$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

$ws.Range("D2").Value = 0.004212758503854275
$ws.Range("E2").Value = 0.0673017748631537
$ws.Range("G2").Value = 0.004141566343605518
$ws.Range("H2").Value = 0.04282997362315655
$ws.Range("I2").Value = 0.001827794127166271
$ws.Range("J2").Value = 0.01286618737503886
$ws.Range("K2").Value = 0.001410416327416897
$ws.Range("D3").Value = 0.001007826998829842
$ws.Range("E3").Value = 0.3842693129554391
$ws.Range("G3").Value = 0.02702527865767479
$ws.Range("H3").Value = 0.2381146480329335
$ws.Range("I3").Value = 0.01914155157282948
$ws.Range("J3").Value = 0.06134908180683851
$ws.Range("K3").Value = 0.009389501065015793
$ws.Range("D4").Value = 0.004408972337841988
$ws.Range("E4").Value = 0.06654203357174993
$ws.Range("G4").Value = 0.004251713398844004
$ws.Range("H4").Value = 0.04214417422190309
$ws.Range("I4").Value = 0.001718732062727213
$ws.Range("J4").Value = 0.01264319103211164
$ws.Range("K4").Value = 0.001474219374358654
$ws.Range("D5").Value = 0.001479329541325569
$ws.Range("E5").Value = 0.383428439963609
$ws.Range("G5").Value = 0.02749304007738829
$ws.Range("H5").Value = 0.2364881676621735
$ws.Range("I5").Value = 0.0199749581515789
$ws.Range("J5").Value = 0.06070742849260569
$ws.Range("K5").Value = 0.009437066502869129
$ws.Range("E6").Value = 1.035875022411346
$ws.Range("D7").Value = 0.004190638661384583
$ws.Range("E7").Value = 0.06552575435489416
$ws.Range("G7").Value = 0.004033647943288088
$ws.Range("H7").Value = 0.04086672281846404
$ws.Range("I7").Value = 0.00188658619299531
$ws.Range("J7").Value = 0.01274234103038907
$ws.Range("K7").Value = 0.001625906210392714
$ws.Range("D8").Value = 0.0009196139872074127
$ws.Range("E8").Value = 0.3518870892003179
$ws.Range("G8").Value = 0.02473711036145687
$ws.Range("H8").Value = 0.2196963313035667
$ws.Range("I8").Value = 0.01750291045755148
$ws.Range("J8").Value = 0.05480776494368911
$ws.Range("K8").Value = 0.008440659381449223
$ws.Range("D9").Value = 0.004432704299688339
$ws.Range("E9").Value = 0.06589768594130874
$ws.Range("G9").Value = 0.004031843040138483
$ws.Range("H9").Value = 0.04140257462859154
$ws.Range("I9").Value = 0.001947587821632624
$ws.Range("J9").Value = 0.01288035791367292
$ws.Range("K9").Value = 0.001401062123477459
$ws.Range("D10").Value = 0.001342307776212692
$ws.Range("E10").Value = 0.3479090337641537
$ws.Range("G10").Value = 0.02436545863747597
$ws.Range("H10").Value = 0.2163986614905298
$ws.Range("I10").Value = 0.01830550003796816
$ws.Range("J10").Value = 0.05400724289938807
$ws.Range("K10").Value = 0.008495531510561705
$ws.Range("E11").Value = 0.9173209187574685
$ws.Range("D12").Value = 0.002281279303133488
$ws.Range("E12").Value = 0.04406110802665353
$ws.Range("G12").Value = 0.002785990480333567
$ws.Range("H12").Value = 0.02779335854575038
$ws.Range("I12").Value = 0.001289741136133671
$ws.Range("J12").Value = 0.008218399249017239
$ws.Range("K12").Value = 0.000989789143204689
$ws.Range("D13").Value = 0.0007802830077707767
$ws.Range("E13").Value = 0.2830429808236659
$ws.Range("G13").Value = 0.01978500094264746
$ws.Range("H13").Value = 0.1763087250292301
$ws.Range("I13").Value = 0.01511721638962626
$ws.Range("J13").Value = 0.04328344948589802
$ws.Range("K13").Value = 0.006938849110156298
$ws.Range("D14").Value = 0.002703116741031408
$ws.Range("E14").Value = 0.05076680891215801
$ws.Range("G14").Value = 0.003207582049071789
$ws.Range("H14").Value = 0.03228841535747051
$ws.Range("I14").Value = 0.001373117789626122
$ws.Range("J14").Value = 0.009126319549977779
$ws.Range("K14").Value = 0.001129476819187403
$ws.Range("D15").Value = 0.001154396682977676
$ws.Range("E15").Value = 0.2906279531307518
$ws.Range("G15").Value = 0.02061355207115412
$ws.Range("H15").Value = 0.1800362728536129
$ws.Range("I15").Value = 0.01625634403899312
$ws.Range("J15").Value = 0.04456278635188937
$ws.Range("K15").Value = 0.007020831573754549
$ws.Range("E16").Value = 1.19781486922875
$ws.Range("D17").Value = 0.002674760762602091
$ws.Range("E17").Value = 0.05211486108601093
$ws.Range("G17").Value = 0.003274162299931049
$ws.Range("H17").Value = 0.03258309187367558
$ws.Range("I17").Value = 0.001271062064915895
$ws.Range("J17").Value = 0.0103395851328969
$ws.Range("K17").Value = 0.001199688762426376
$ws.Range("D18").Value = 0.0008676475845277309
$ws.Range("E18").Value = 0.3214176730252802
$ws.Range("G18").Value = 0.0224135834723711
$ws.Range("H18").Value = 0.1991616445593536
$ws.Range("I18").Value = 0.01742002135142684
$ws.Range("J18").Value = 0.05026417504996061
$ws.Range("K18").Value = 0.007777146995067596
$ws.Range("D19").Value = 0.002820469439029694
$ws.Range("E19").Value = 0.0530701931566
$ws.Range("G19").Value = 0.003372336272150278
$ws.Range("H19").Value = 0.03259953297674656
$ws.Range("I19").Value = 0.001340185292065144
$ws.Range("J19").Value = 0.01092933863401413
$ws.Range("K19").Value = 0.001208415254950523
$ws.Range("D20").Value = 0.001317867077887058
$ws.Range("E20").Value = 0.3277451996691525
$ws.Range("G20").Value = 0.02273847255855799
$ws.Range("H20").Value = 0.2024244735948741
$ws.Range("I20").Value = 0.01867597969248891
$ws.Range("J20").Value = 0.05149086331948638
$ws.Range("K20").Value = 0.007962895557284355
$ws.Range("E21").Value = 1.015661107841879
$ws.Range("D22").Value = 0.002665614243596792
$ws.Range("E22").Value = 0.04612505994737148
$ws.Range("G22").Value = 0.002748556435108185
$ws.Range("H22").Value = 0.02920609433203936
$ws.Range("I22").Value = 0.00117412069812417
$ws.Range("J22").Value = 0.009055460337549448
$ws.Range("K22").Value = 0.001013088040053844
$ws.Range("D23").Value = 0.0007669990882277489
$ws.Range("E23").Value = 0.2942635361105204
$ws.Range("G23").Value = 0.02059155749157071
$ws.Range("H23").Value = 0.1811021571047604
$ws.Range("I23").Value = 0.0163769512437284
$ws.Range("J23").Value = 0.04608120024204254
$ws.Range("K23").Value = 0.007416658569127321
$ws.Range("D24").Value = 0.002863870933651924
$ws.Range("E24").Value = 0.04482023511081934
$ws.Range("G24").Value = 0.002723724581301212
$ws.Range("H24").Value = 0.02809604303911328
$ws.Range("I24").Value = 0.001164416316896677
$ws.Range("J24").Value = 0.008978911675512791
$ws.Range("K24").Value = 0.0009899823926389217
$ws.Range("D25").Value = 0.001167423091828823
$ws.Range("E25").Value = 0.3058887221850455
$ws.Range("G25").Value = 0.02149605099111795
$ws.Range("H25").Value = 0.1878059906885028
$ws.Range("I25").Value = 0.0177012630738318
$ws.Range("J25").Value = 0.04820545297116041
$ws.Range("K25").Value = 0.007510603405535221
$ws.Range("E26").Value = 1.154503941070288
$ws.Range("D27").Value = 0.003205010294914246
$ws.Range("E27").Value = 0.05502571165561676
$ws.Range("G27").Value = 0.003396784700453281
$ws.Range("H27").Value = 0.03465584823861718
$ws.Range("I27").Value = 0.001489860843867063
$ws.Range("J27").Value = 0.0106443946249783
$ws.Range("K27").Value = 0.001247777696698904
$ws.Range("D28").Value = 0.0008684741333127022
$ws.Range("E28").Value = 0.3269761184230447
$ws.Range("G28").Value = 0.02291050618514419
$ws.Range("H28").Value = 0.2028767012059688
$ws.Range("I28").Value = 0.01711173020303249
$ws.Range("J28").Value = 0.05115713430568576
$ws.Range("K28").Value = 0.007992563024163245
$ws.Range("D29").Value = 0.00344582675024867
$ws.Range("E29").Value = 0.0562193913385272
$ws.Range("G29").Value = 0.003517439868301153
$ws.Range("H29").Value = 0.03530614804476499
$ws.Range("I29").Value = 0.001508807856589556
$ws.Range("J29").Value = 0.01091162376105785
$ws.Range("K29").Value = 0.001240631192922592
$ws.Range("D30").Value = 0.001292264834046364
$ws.Range("E30").Value = 0.3311198697425425
$ws.Range("G30").Value = 0.02334131486713886
$ws.Range("H30").Value = 0.2046307132579386
$ws.Range("I30").Value = 0.01818280899897218
$ws.Range("J30").Value = 0.05179475480690598
$ws.Range("K30").Value = 0.008085385710000993
$ws.Range("E31").Value = 1.064235171861946
